$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2088.9375
$ws.Range("I43").Value = 2243.4167
$ws.Range("J43").Value = 1625.5
$ws.Range("K43").Value = 2243.4167
$ws.Range("L43").Value = 1625.5
$ws.Range("M43").Value = -2174.4167
$ws.Range("N43").Value = -1763.5

$ws.Range("H98").Value = 598
$ws.Range("I98").Value = 559
$ws.Range("K98").Value = 559
$ws.Range("M98").Value = 939

$ws.Range("H106").Value = 4154.857
$ws.Range("I106").Value = 4451
$ws.Range("J106").Value = 3760
$ws.Range("K106").Value = 4451
$ws.Range("L106").Value = 3760
$ws.Range("M106").Value = -3820
$ws.Range("N106").Value = -5022

$ws.Range("H111").Value = 5558721.5
$ws.Range("I111").Value = 3477.6155
$ws.Range("J111").Value = 20002356
$ws.Range("K111").Value = 10432.8465
$ws.Range("L111").Value = 60007068
$ws.Range("M111").Value = -7365.8465
$ws.Range("N111").Value = -60013202

$ws.Range("H112").Value = 1214.2979
$ws.Range("J112").Value = 1229.826
$ws.Range("L112").Value = 3689.478
$ws.Range("N112").Value = -5905.478

$ws.Range("H113").Value = 113222.78
$ws.Range("I113").Value = 145000.72
$ws.Range("K113").Value = 145000.72
$ws.Range("M113").Value = -141746.72

$ws.Range("H122").Value = 598
$ws.Range("I122").Value = 559
$ws.Range("K122").Value = 1677
$ws.Range("M122").Value = 773

$ws.Range("H133").Value = 51991.668
$ws.Range("J133").Value = 51991.668
$ws.Range("L133").Value = 51991.668
$ws.Range("N133").Value = -62111.668

$ws.Range("H135").Value = 1370.4166
$ws.Range("I135").Value = 575.59375
$ws.Range("J135").Value = 2960.0625
$ws.Range("K135").Value = 5180.34375
$ws.Range("L135").Value = 26640.5625
$ws.Range("M135").Value = -2645.34375
$ws.Range("N135").Value = -31710.5625

$ws.Range("H137").Value = 2066.3333
$ws.Range("I137").Value = 1843.25
$ws.Range("K137").Value = 5529.75
$ws.Range("M137").Value = -2979.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 144927.28
$ws.Range("I2").Value = 2298.2
$ws.Range("J2").Value = 501500
$ws.Range("K2").Value = 2298.2
$ws.Range("L2").Value = 501500
$ws.Range("M2").Value = -2185.2
$ws.Range("N2").Value = -501726

$ws.Range("H61").Value = 1461.16
$ws.Range("I61").Value = 1143.5714
$ws.Range("K61").Value = 1143.5714
$ws.Range("M61").Value = -931.5714

$ws.Range("H116").Value = 144927.28
$ws.Range("I116").Value = 2298.2
$ws.Range("J116").Value = 501500
$ws.Range("K116").Value = 2298.2
$ws.Range("L116").Value = 501500
$ws.Range("M116").Value = -4.199999999999818
$ws.Range("N116").Value = -506088

$ws.Range("H132").Value = 15033.682
$ws.Range("I132").Value = 16189.077
$ws.Range("J132").Value = 6021.6
$ws.Range("K132").Value = 48567.231
$ws.Range("L132").Value = 18064.8
$ws.Range("M132").Value = -46037.231
$ws.Range("N132").Value = -23124.8

$ws.Range("H136").Value = 1461.16
$ws.Range("I136").Value = 1143.5714
$ws.Range("K136").Value = 3430.7142
$ws.Range("M136").Value = -880.7142000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 144927.28
$ws.Range("I3").Value = 2298.2
$ws.Range("J3").Value = 501500
$ws.Range("K3").Value = 2298.2
$ws.Range("L3").Value = 501500
$ws.Range("M3").Value = -2184.2
$ws.Range("N3").Value = -501728

$ws.Range("H107").Value = 47620044
$ws.Range("I107").Value = 62501000
$ws.Range("J107").Value = 984.4
$ws.Range("K107").Value = 62501000
$ws.Range("L107").Value = 984.4
$ws.Range("M107").Value = -62499080
$ws.Range("N107").Value = -4824.4

$ws.Range("H134").Value = 10901.193
$ws.Range("I134").Value = 11449.52
$ws.Range("J134").Value = 5198.6
$ws.Range("K134").Value = 34348.56
$ws.Range("L134").Value = 15595.8
$ws.Range("M134").Value = -31813.56
$ws.Range("N134").Value = -20665.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 67314.87
$ws.Range("I16").Value = 615.63635
$ws.Range("K16").Value = 615.63635
$ws.Range("M16").Value = -328.63635

$ws.Range("H20").Value = 48675.6
$ws.Range("J20").Value = 48675.6
$ws.Range("L20").Value = 48675.6
$ws.Range("N20").Value = -49147.6

$ws.Range("H30").Value = 48675.6
$ws.Range("J30").Value = 48675.6
$ws.Range("L30").Value = 48675.6
$ws.Range("N30").Value = -48857.6

$ws.Range("H107").Value = 832.5333000000001
$ws.Range("I107").Value = 932.7895
$ws.Range("J107").Value = 659.36365
$ws.Range("K107").Value = 932.7895
$ws.Range("L107").Value = 659.36365
$ws.Range("M107").Value = 987.2105
$ws.Range("N107").Value = -4499.36365

$ws.Range("H113").Value = 67314.87
$ws.Range("I113").Value = 615.63635
$ws.Range("K113").Value = 615.63635
$ws.Range("M113").Value = 1554.36365

$ws.Range("H128").Value = 48675.6
$ws.Range("J128").Value = 48675.6
$ws.Range("L128").Value = 48675.6
$ws.Range("N128").Value = -58635.6

$ws.Range("H132").Value = 62503084
$ws.Range("I132").Value = 76926790
$ws.Range("J132").Value = 45456884
$ws.Range("K132").Value = 230780370
$ws.Range("L132").Value = 136370652
$ws.Range("M132").Value = -230777840
$ws.Range("N132").Value = -136375712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 811936.25
$ws.Range("I131").Value = 538.25
$ws.Range("J131").Value = 1002853.4
$ws.Range("K131").Value = 1614.75
$ws.Range("L131").Value = 3008560.2
$ws.Range("M131").Value = 3425.25
$ws.Range("N131").Value = -3018640.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2250
$ws.Range("I113").Value = 2250
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -80
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 689761.8
$ws.Range("I16").Value = 144314.58
$ws.Range("J16").Value = 1113998.5
$ws.Range("K16").Value = 144314.58
$ws.Range("L16").Value = 1113998.5
$ws.Range("M16").Value = -144144.58
$ws.Range("N16").Value = -1114338.5

$ws.Range("H61").Value = 2637.3635
$ws.Range("I61").Value = 2667.889
$ws.Range("K61").Value = 2667.889
$ws.Range("M61").Value = -2465.889

$ws.Range("H113").Value = 2637.3635
$ws.Range("I113").Value = 2667.889
$ws.Range("K113").Value = 2667.889
$ws.Range("M113").Value = -497.8890000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1578.6666
$ws.Range("I12").Value = 736
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 736
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -594
$ws.Range("N12").Value = -2284

$ws.Range("H100").Value = 73118.14
$ws.Range("I100").Value = 101535.4
$ws.Range("J100").Value = 2075
$ws.Range("K100").Value = 203070.8
$ws.Range("L100").Value = 4150
$ws.Range("M100").Value = -202529.8
$ws.Range("N100").Value = -5232

$ws.Range("H113").Value = 705.7273
$ws.Range("I113").Value = 532
$ws.Range("J113").Value = 850.5
$ws.Range("K113").Value = 1596
$ws.Range("L113").Value = 2551.5
$ws.Range("M113").Value = 574
$ws.Range("N113").Value = -6891.5

$ws.Range("H119").Value = 27607.445
$ws.Range("J119").Value = 27607.445
$ws.Range("L119").Value = 27607.445
$ws.Range("N119").Value = -37283.445

$ws.Range("H122").Value = 2598.9546
$ws.Range("I122").Value = 1954
$ws.Range("K122").Value = 5862
$ws.Range("M122").Value = -3412

$ws.Range("H136").Value = 15361.025
$ws.Range("I136").Value = 27370.81
$ws.Range("J136").Value = 4522.927
$ws.Range("K136").Value = 82112.43000000001
$ws.Range("L136").Value = 13568.781
$ws.Range("M136").Value = -79562.43000000001
$ws.Range("N136").Value = -18668.781
